# Daily attendance processing - 2025-10-18 03:29:20
# Swap the order of "<name>, System" -> "System, <name>" in the
# "Recorded By" column (G) of the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 1; $i -le $lastRow; $i++) {
    $cell = $ws.Range("G$i")
    $val = $cell.Value2

    if ($val -ne $null -and $val -match "^[^,]+, System$") {
        $newVal = $val -replace "^(.+), System$", "System, `$1"
        $cell.Value = $newVal
    }
}
